$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ row=2; A="ECs"; B="Pdgfb"; C="Pdgfra"; D="ECs"; E=2; F=0.6666666666666666; G=27.388319; H=82.164957; I=0.9327824776802173; J=0.9327824776802174; K=2; L=0.6666666666666666; M=1.155747666666667; N=3.467243; O=0.004246591903937912; P=0.004246591903937912; Q=31.65398577817233; R=284.885872003551; S=0.003961146517851957; T=0.003961146517851958 },
  @{ row=3; A="ECs"; B="Pdgfb"; C="Pdgfra"; D="FAPs"; E=2; F=0.6666666666666666; G=27.388319; H=82.164957; I=0.9327824776802173; J=0.9327824776802174; K=3; L=1; M=270.7963256666667; N=812.3889770000001; O=0.9949935590256014; P=0.9949935590256014; Q=7416.656151386554; R=66749.905362479; S=0.928112557263758; T=0.9281125572637581 },
  @{ row=4; A="ECs"; B="Pdgfb"; C="Pdgfra"; D="sCs"; E=2; F=0.6666666666666666; G=27.388319; H=82.164957; I=0.9327824776802173; J=0.9327824776802174; K=3; L=1; M=0.2067996666666667; N=0.620399; O=0.0007598490704606447; P=0.0007598490704606446; Q=5.663895239760333; R=50.975057157843; S=0.0007087738986072902; T=0.0007087738986072901 },
  @{ row=5; A="FAPs"; B="Pdgfb"; C="Pdgfra"; D="ECs"; E=2; F=0.6666666666666666; G=0.134644; H=0.403932; I=0.004585661643738528; J=0.004585661643738528; K=2; L=0.6666666666666666; M=1.155747666666667; N=3.467243; O=0.004246591903937912; P=0.004246591903937912; Q=0.1556144888306666; R=1.400530399476; S=[double]"1.947343361049865E-05"; T=[double]"1.947343361049865E-05" },
  @{ row=6; A="FAPs"; B="Pdgfb"; C="Pdgfra"; D="FAPs"; E=2; F=0.6666666666666666; G=0.134644; H=0.403932; I=0.004585661643738528; J=0.004585661643738528; K=3; L=1; M=270.7963256666667; N=812.3889770000001; O=0.9949935590256014; P=0.9949935590256014; Q=36.46110047306266; R=328.149904257564; S=0.004562703799390587; T=0.004562703799390587 },
  @{ row=7; A="FAPs"; B="Pdgfb"; C="Pdgfra"; D="sCs"; E=2; F=0.6666666666666666; G=0.134644; H=0.403932; I=0.004585661643738528; J=0.004585661643738528; K=3; L=1; M=0.2067996666666667; N=0.620399; O=0.0007598490704606447; P=0.0007598490704606446; Q=0.02784433431866667; R=0.250599008868; S=[double]"3.484410737441753E-06"; T=[double]"3.484410737441753E-06" },
  @{ row=8; A="sCs"; B="Pdgfb"; C="Pdgfra"; D="ECs"; E=3; F=1; G=1.838994; H=5.516982; I=0.06263186067604418; J=0.06263186067604418; K=2; L=0.6666666666666666; M=1.155747666666667; N=3.467243; O=0.004246591903937912; P=0.004246591903937912; Q=2.125413024514; R=19.128717220626; S=0.0002659719524754565; T=0.0002659719524754565 },
  @{ row=9; A="sCs"; B="Pdgfb"; C="Pdgfra"; D="FAPs"; E=3; F=1; G=1.838994; H=5.516982; I=0.06263186067604418; J=0.06263186067604418; K=3; L=1; M=270.7963256666667; N=812.3889770000001; O=0.9949935590256014; P=0.9949935590256014; Q=497.9928181230461; R=4481.935363107415; S=0.06231829796245281; T=0.06231829796245281 },
  @{ row=10; A="sCs"; B="Pdgfb"; C="Pdgfra"; D="sCs"; E=3; F=1; G=1.838994; H=5.516982; I=0.06263186067604418; J=0.06263186067604418; K=3; L=1; M=0.2067996666666667; N=0.620399; O=0.0007598490704606447; P=0.0007598490704606446; Q=0.3803033462020001; R=3.422730115818001; S=[double]"4.759076111591278E-05"; T=[double]"4.759076111591277E-05" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.row, 1).Value = $r.A
    $ws.Cells.Item($r.row, 2).Value = $r.B
    $ws.Cells.Item($r.row, 3).Value = $r.C
    $ws.Cells.Item($r.row, 4).Value = $r.D
    $ws.Cells.Item($r.row, 5).Value = $r.E
    $ws.Cells.Item($r.row, 6).Value = $r.F
    $ws.Cells.Item($r.row, 7).Value = $r.G
    $ws.Cells.Item($r.row, 8).Value = $r.H
    $ws.Cells.Item($r.row, 9).Value = $r.I
    $ws.Cells.Item($r.row, 10).Value = $r.J
    $ws.Cells.Item($r.row, 11).Value = $r.K
    $ws.Cells.Item($r.row, 12).Value = $r.L
    $ws.Cells.Item($r.row, 13).Value = $r.M
    $ws.Cells.Item($r.row, 14).Value = $r.N
    $ws.Cells.Item($r.row, 15).Value = $r.O
    $ws.Cells.Item($r.row, 16).Value = $r.P
    $ws.Cells.Item($r.row, 17).Value = $r.Q
    $ws.Cells.Item($r.row, 18).Value = $r.R
    $ws.Cells.Item($r.row, 19).Value = $r.S
    $ws.Cells.Item($r.row, 20).Value = $r.T
}
